$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2021-11-15"

$ws.Range("A13").Value = "November (through 11-15)"

$ws.Range("H13").Value = 2
$ws.Range("I13").Value = 62
$ws.Range("J13").Value = 0.0312
$ws.Range("L13").Value = 27
$ws.Range("M13").Value = 0.1562
$ws.Range("O13").Value = 21
$ws.Range("P13").Value = 0.16
$ws.Range("R13").Value = 85
$ws.Range("S13").Value = 0.0341
$ws.Range("U13").Value = 107
$ws.Range("V13").Value = 0.0093

$ws.Range("H14").Value = 63
$ws.Range("I14").Value = 711
$ws.Range("J14").Value = 0.0814
$ws.Range("L14").Value = 576
$ws.Range("M14").Value = 0.1097
$ws.Range("O14").Value = 455
$ws.Range("P14").Value = 0.1026
$ws.Range("R14").Value = 1088
$ws.Range("U14").Value = 1461
$ws.Range("V14").Value = 0.0574
